# =====================================================================
# Edit script: adds player-info / "ODI Batting Extra" sheets, renames
# MATCH_CARD_LINK -> MATCH_CODE on the two pre-existing sheets and
# replaces the full scorecard URLs with the bare match codes.
# =====================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert a new "Player Info" sheet in front of everything else
# ---------------------------------------------------------------------
$firstSheet = $wb.Worksheets.Item(1)
$wsInfo = $wb.Worksheets.Add($firstSheet)
$wsInfo.Name = "Player Info"

$infoHeaders = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($i = 0; $i -lt $infoHeaders.Length; $i++) {
    $cell = $wsInfo.Cells.Item(1, $i + 1)
    $cell.Value = $infoHeaders[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

$wsInfo.Cells.Item(2, 1).Value = "'4789"
$wsInfo.Cells.Item(2, 2).Value = "Hamza Tahir"
$wsInfo.Cells.Item(2, 3).Value = "Right Handed"
$wsInfo.Cells.Item(2, 4).Value = "Left Arm Orthodox"

# ---------------------------------------------------------------------
# 2. "ODI Batting" sheet: rename MATCH_CARD_LINK -> MATCH_CODE and
#    replace the match-card URLs with the bare numeric match code
# ---------------------------------------------------------------------
$wsBatting = $wb.Worksheets.Item("ODI Batting")
$wsBatting.Cells.Item(1, 4).Value = "MATCH_CODE"

$battingCodes = @{
    2  = "4364"
    3  = "4365"
    4  = "4366"
    5  = "4381"
    6  = "4384"
    7  = "4386"
    8  = "4461"
    9  = "4510"
    10 = "4512"
    11 = "4513"
    12 = "4515"
    13 = "4569"
    14 = "4570"
    15 = "4572"
    16 = "4573"
    17 = "4575"
    18 = "4576"
    19 = "4578"
    20 = "4604"
    21 = "4610"
    22 = "4612"
    23 = "4617"
    24 = "4625"
    25 = "4629"
    26 = "4632"
    27 = "4635"
    28 = "4677"
    29 = "4681"
    30 = "4680"
    31 = "4702"
    32 = "4703"
}

foreach ($r in $battingCodes.Keys) {
    $wsBatting.Cells.Item($r, 4).Value = "'" + $battingCodes[$r]
}

# ---------------------------------------------------------------------
# 3. "ODI Bowling" sheet: rename MATCH_CARD_LINK -> MATCH_CODE and
#    replace the match-card URLs with the bare numeric match code
# ---------------------------------------------------------------------
$wsBowling = $wb.Worksheets.Item("ODI Bowling")
$wsBowling.Cells.Item(1, 2).Value = "MATCH_CODE"

$bowlingCodes = @{
    2  = "4364"
    3  = "4365"
    4  = "4366"
    5  = "4381"
    6  = "4384"
    7  = "4386"
    8  = "4461"
    9  = "4510"
    10 = "4512"
    11 = "4513"
    12 = "4569"
    13 = "4570"
    14 = "4572"
    15 = "4573"
    16 = "4575"
    17 = "4576"
    18 = "4578"
    19 = "4604"
    20 = "4610"
    21 = "4612"
    22 = "4617"
    23 = "4625"
    24 = "4629"
    25 = "4632"
    26 = "4635"
    27 = "4677"
    28 = "4681"
    29 = "4680"
    30 = "4702"
    31 = "4703"
}

foreach ($r in $bowlingCodes.Keys) {
    $wsBowling.Cells.Item($r, 2).Value = "'" + $bowlingCodes[$r]
}

# ---------------------------------------------------------------------
# 4. Append a new "ODI Batting Extra" sheet after "ODI Bowling"
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsExtra = $wb.Worksheets.Add($null, $lastSheet)
$wsExtra.Name = "ODI Batting Extra"

$extraHeaders = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($i = 0; $i -lt $extraHeaders.Length; $i++) {
    $cell = $wsExtra.Cells.Item(1, $i + 1)
    $cell.Value = $extraHeaders[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# match_code, batting_position, num4, num6, percent_runs_of_total, man_of_match
$extraData = @(
    @("4569", "8",  "",  "",  "",       "NO"),
    @("4570", "11", "0", "0", "",       "NO"),
    @("4572", "11", "",  "",  "",       "NO"),
    @("4573", "",   "",  "",  "",       "NO"),
    @("4575", "11", "",  "",  "",       "NO"),
    @("4576", "10", "",  "",  "",       "NO"),
    @("4578", "11", "",  "",  "",       "NO"),
    @("4604", "11", "0", "0", "4.86%",  "NO"),
    @("4610", "11", "1", "0", "1.63%",  "NO"),
    @("4612", "11", "0", "0", "",       "NO"),
    @("4617", "11", "",  "",  "",       "NO"),
    @("4625", "",   "",  "",  "",       "NO"),
    @("4629", "11", "0", "0", "0.48%",  "NO"),
    @("4632", "11", "",  "",  "",       "NO"),
    @("4635", "",   "",  "",  "",       ""),
    @("4677", "",   "",  "",  "",       ""),
    @("4681", "",   "",  "",  "",       ""),
    @("4680", "",   "",  "",  "",       ""),
    @("4702", "",   "",  "",  "",       ""),
    @("4703", "",   "",  "",  "",       "")
)

for ($i = 0; $i -lt $extraData.Length; $i++) {
    $r = $i + 2
    $row = $extraData[$i]

    $wsExtra.Cells.Item($r, 1).Value = "'" + $row[0]

    if ($row[1] -ne "") {
        $wsExtra.Cells.Item($r, 2).Value = [double]$row[1]
    }

    if ($row[2] -ne "") {
        $wsExtra.Cells.Item($r, 3).Value = "'" + $row[2]
    }

    if ($row[3] -ne "") {
        $wsExtra.Cells.Item($r, 4).Value = "'" + $row[3]
    }

    if ($row[4] -ne "") {
        $wsExtra.Cells.Item($r, 5).Value = "'" + $row[4]
    }

    if ($row[5] -ne "") {
        $wsExtra.Cells.Item($r, 6).Value = $row[5]
    }
}

# ---------------------------------------------------------------------
# 5. Keep the first sheet active, matching the original workbook state
# ---------------------------------------------------------------------
$wsInfo.Activate()
